$d = $word.ActiveDocument

$p1 = $d.Paragraphs(1).Range

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
       '<w:body>' +
       '<w:p w14:paraId="5ADF5830" w14:textId="42E3A3E7" w:rsidR="00384372" w:rsidRDefault="00094D0B">' +
       '<w:r><w:t>This is a Microsoft word document.</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' +
       '<w:r><w:t>Changed main</w:t></w:r>' +
       '<w:r><w:t>)</w:t></w:r>' +
       '</w:p>' +
       '</w:body></w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$p1.InsertXML($xml)
